# Update "想去人数" (want-to-go count, column F) figures across sheets,
# matching the data refresh recorded in the commit
# "Update gh-pages to output generated at 456a3b4".
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 923
$ws.Range("F3").Value = 1026
$ws.Range("F4").Value = 811
$ws.Range("F5").Value = 887
$ws.Range("F6").Value = 468
$ws.Range("F7").Value = 717
$ws.Range("F8").Value = 168
$ws.Range("F9").Value = 1320
$ws.Range("F10").Value = 744
$ws.Range("F11").Value = 426
$ws.Range("F12").Value = 567
$ws.Range("F14").Value = 58
$ws.Range("F15").Value = 1209
$ws.Range("F16").Value = 147
$ws.Range("F17").Value = 20
$ws.Range("F18").Value = 431
$ws.Range("F19").Value = 380
$ws.Range("F20").Value = 97
$ws.Range("F21").Value = 606
$ws.Range("F22").Value = 162
$ws.Range("F24").Value = 38
$ws.Range("F25").Value = 1097
$ws.Range("F26").Value = 18

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 259
$ws.Range("F9").Value = 31
$ws.Range("F11").Value = 117

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 923
$ws.Range("F5").Value = 1026
$ws.Range("F6").Value = 811
$ws.Range("F7").Value = 887
$ws.Range("F8").Value = 468
$ws.Range("F9").Value = 468
$ws.Range("F10").Value = 717
$ws.Range("F11").Value = 168
$ws.Range("F12").Value = 1320
$ws.Range("F13").Value = 744
$ws.Range("F16").Value = 426
$ws.Range("F17").Value = 567
$ws.Range("F20").Value = 58
$ws.Range("F21").Value = 1209
$ws.Range("F23").Value = 147
$ws.Range("F24").Value = 20
$ws.Range("F25").Value = 431
$ws.Range("F26").Value = 380
$ws.Range("F27").Value = 97
$ws.Range("F28").Value = 259
$ws.Range("F30").Value = 606
$ws.Range("F31").Value = 31
$ws.Range("F33").Value = 117
$ws.Range("F34").Value = 117
$ws.Range("F35").Value = 162
$ws.Range("F37").Value = 38
$ws.Range("F38").Value = 1097
$ws.Range("F39").Value = 18
